# Update countries & provincias Spain
# Refreshed COVID-19 case counts for several countries; since the sheet is
# kept sorted descending by "Casos totales" (column B), a handful of
# countries changed rank and the rows that hold them were overwritten with
# the country that now occupies that rank, plus the new case numbers.
# The "last updated" footer timestamp is also bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp (row 1, column A)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 6 de Abril de 2020 a las 06:22"

# Row 10: Iran (same rank, updated case numbers)
$ws.Cells.Item(10, 1).Value = "Iran"
$ws.Cells.Item(10, 2).Value = 58226
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 22011
$ws.Cells.Item(10, 5).Value = 32612
$ws.Cells.Item(10, 6).Value = 4103
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 3603

# Row 23: Australia (same rank, updated case numbers)
$ws.Cells.Item(23, 1).Value = "Australia"
$ws.Cells.Item(23, 2).Value = 5773
$ws.Cells.Item(23, 3).Value = 23
$ws.Cells.Item(23, 4).Value = 2315
$ws.Cells.Item(23, 5).Value = 3419
$ws.Cells.Item(23, 6).Value = 91
$ws.Cells.Item(23, 7).Value = 2
$ws.Cells.Item(23, 8).Value = 39

# Row 24: Noruega (same rank, updated case numbers)
$ws.Cells.Item(24, 1).Value = "Noruega"
$ws.Cells.Item(24, 2).Value = 5759
$ws.Cells.Item(24, 3).Value = 72
$ws.Cells.Item(24, 4).Value = 32
$ws.Cells.Item(24, 5).Value = 5656
$ws.Cells.Item(24, 6).Value = 89
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 71

# Row 30: India (same rank, updated case numbers)
$ws.Cells.Item(30, 1).Value = "India"
$ws.Cells.Item(30, 2).Value = 4289
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 328
$ws.Cells.Item(30, 5).Value = 3843
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 118

# Row 36: Pakistan moves up one rank (was Filipinas)
$ws.Cells.Item(36, 1).Value = "Pakistan"
$ws.Cells.Item(36, 2).Value = 3277
$ws.Cells.Item(36, 3).Value = 120
$ws.Cells.Item(36, 4).Value = 257
$ws.Cells.Item(36, 5).Value = 2970
$ws.Cells.Item(36, 6).Value = 18
$ws.Cells.Item(36, 7).Value = 3
$ws.Cells.Item(36, 8).Value = 50

# Row 37: Filipinas moves down one rank (was Pakistan)
$ws.Cells.Item(37, 1).Value = "Filipinas"
$ws.Cells.Item(37, 2).Value = 3246
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = 64
$ws.Cells.Item(37, 5).Value = 3030
$ws.Cells.Item(37, 6).Value = 1
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 152

# Row 88: Afganistan moves up two ranks (was Taiwan)
$ws.Cells.Item(88, 1).Value = "Afganistan"
$ws.Cells.Item(88, 2).Value = 367
$ws.Cells.Item(88, 3).Value = 18
$ws.Cells.Item(88, 4).Value = 15
$ws.Cells.Item(88, 5).Value = 345
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 7

# Row 89: Taiwan moves down one rank (was Albania)
$ws.Cells.Item(89, 1).Value = "Taiwan"
$ws.Cells.Item(89, 2).Value = 363
$ws.Cells.Item(89, 3).Value = 0
$ws.Cells.Item(89, 4).Value = 54
$ws.Cells.Item(89, 5).Value = 304
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 5

# Row 90: Albania moves down one rank (was Afganistan)
$ws.Cells.Item(90, 1).Value = "Albania"
$ws.Cells.Item(90, 2).Value = 361
$ws.Cells.Item(90, 3).Value = 0
$ws.Cells.Item(90, 4).Value = 104
$ws.Cells.Item(90, 5).Value = 237
$ws.Cells.Item(90, 6).Value = 7
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 20

# Row 106: Kirguistan moves up several ranks (was Montenegro)
$ws.Cells.Item(106, 1).Value = "Kirguistan"
$ws.Cells.Item(106, 2).Value = 216
$ws.Cells.Item(106, 3).Value = 69
$ws.Cells.Item(106, 4).Value = 9
$ws.Cells.Item(106, 5).Value = 206
$ws.Cells.Item(106, 6).Value = 5
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 1

# Row 107: Montenegro moves down one rank (was Ghana)
$ws.Cells.Item(107, 1).Value = "Montenegro"
$ws.Cells.Item(107, 2).Value = 214
$ws.Cells.Item(107, 3).Value = 0
$ws.Cells.Item(107, 4).Value = 1
$ws.Cells.Item(107, 5).Value = 211
$ws.Cells.Item(107, 6).Value = 4
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 2

# Row 108: Ghana moves down one rank (was Niger)
$ws.Cells.Item(108, 1).Value = "Ghana"
$ws.Cells.Item(108, 2).Value = 214
$ws.Cells.Item(108, 3).Value = 0
$ws.Cells.Item(108, 4).Value = 31
$ws.Cells.Item(108, 5).Value = 178
$ws.Cells.Item(108, 6).Value = 2
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 5

# Row 109: Niger moves down one rank (was Bolivia)
$ws.Cells.Item(109, 1).Value = "Niger"
$ws.Cells.Item(109, 2).Value = 184
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(109, 4).Value = 13
$ws.Cells.Item(109, 5).Value = 161
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 10

# Row 110: Bolivia moves down one rank (was Islas Feroe)
$ws.Cells.Item(110, 1).Value = "Bolivia"
$ws.Cells.Item(110, 2).Value = 183
$ws.Cells.Item(110, 3).Value = 26
$ws.Cells.Item(110, 4).Value = 2
$ws.Cells.Item(110, 5).Value = 170
$ws.Cells.Item(110, 6).Value = 3
$ws.Cells.Item(110, 7).Value = 1
$ws.Cells.Item(110, 8).Value = 11

# Row 111: Islas Feroe moves down one rank (was Sri Lanka)
$ws.Cells.Item(111, 1).Value = "Islas Feroe"
$ws.Cells.Item(111, 2).Value = 181
$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(111, 4).Value = 99
$ws.Cells.Item(111, 5).Value = 82
$ws.Cells.Item(111, 6).Value = 1
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 0

# Row 112: Sri Lanka moves down one rank (was Georgia)
$ws.Cells.Item(112, 1).Value = "Sri Lanka"
$ws.Cells.Item(112, 2).Value = 176
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(112, 4).Value = 33
$ws.Cells.Item(112, 5).Value = 138
$ws.Cells.Item(112, 6).Value = 5
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 5

# Row 113: Georgia moves down one rank (was Venezuela)
$ws.Cells.Item(113, 1).Value = "Georgia"
$ws.Cells.Item(113, 2).Value = 174
$ws.Cells.Item(113, 3).Value = 0
$ws.Cells.Item(113, 4).Value = 36
$ws.Cells.Item(113, 5).Value = 136
$ws.Cells.Item(113, 6).Value = 6
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 2

# Row 114: Venezuela moves down one rank (was Consejo Danes para los Refugiados)
$ws.Cells.Item(114, 1).Value = "Venezuela"
$ws.Cells.Item(114, 2).Value = 159
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 52
$ws.Cells.Item(114, 5).Value = 100
$ws.Cells.Item(114, 6).Value = 6
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 7

# Row 115: Consejo Danes para los Refugiados moves down one rank (was Martinica)
$ws.Cells.Item(115, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(115, 2).Value = 154
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 4).Value = 3
$ws.Cells.Item(115, 5).Value = 133
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 18

# Row 116: Martinica moves down one rank (was Kirguistan)
$ws.Cells.Item(116, 1).Value = "Martinica"
$ws.Cells.Item(116, 2).Value = 149
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 50
$ws.Cells.Item(116, 5).Value = 95
$ws.Cells.Item(116, 6).Value = 21
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 4

# Row 133: El Salvador moves up two ranks (was Guayana Francesa)
$ws.Cells.Item(133, 1).Value = "El Salvador"
$ws.Cells.Item(133, 2).Value = 69
$ws.Cells.Item(133, 3).Value = 7
$ws.Cells.Item(133, 4).Value = 2
$ws.Cells.Item(133, 5).Value = 64
$ws.Cells.Item(133, 6).Value = 4
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 3

# Row 134: Guayana Francesa moves down one rank (was Aruba)
$ws.Cells.Item(134, 1).Value = "Guayana Francesa"
$ws.Cells.Item(134, 2).Value = 68
$ws.Cells.Item(134, 3).Value = 0
$ws.Cells.Item(134, 4).Value = 27
$ws.Cells.Item(134, 5).Value = 41
$ws.Cells.Item(134, 6).Value = 1
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 0

# Row 135: Aruba moves down one rank (was El Salvador)
$ws.Cells.Item(135, 1).Value = "Aruba"
$ws.Cells.Item(135, 2).Value = 64
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(135, 4).Value = 1
$ws.Cells.Item(135, 5).Value = 63
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 0
